# Update Jogos_do_Dia_Betfair_Back_Lay_2025-12-19 worksheet:
#  - insert a new match row (Azerbaijan Premier League) at row 4, shifting the
#    remaining matches down by one row
#  - tweak a handful of odds that changed on the surrounding rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new match by inserting a row at position 4.
#    This shifts the old rows 4-9 down to 5-10 (and grows the used range to
#    A1:AO10), matching the sheet-level "dimension" change in the diff.
$ws.Rows.Item(4).Insert()

# 2) Populate the newly inserted row 4 with the Azerbaijan Premier League match.
#    The date column uses a leading apostrophe so Excel stores "2025-12-19" as
#    plain text instead of auto-converting it to a real date serial (matching
#    every other row in the column); the Style reset afterwards drops the
#    transient "quote prefix" formatting so the cell stays plain/unstyled.
$ws.Range("A4").Value = "Azerbaijan Premier League"
$ws.Range("B4").Value = "'2025-12-19"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "07:00:00"
$ws.Range("D4").Value = "Kapaz Ganja"
$ws.Range("E4").Value = "PFK Turan Tovuz"
$ws.Range("F4").Value = 1.04
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 1.04
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1.02
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.3
$ws.Range("O4").Value = 1.01
$ws.Range("P4").Value = 1.3
$ws.Range("Q4").Value = 1.32
$ws.Range("R4").Value = 1.18
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 1.01
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# 3) A few odds on the still-untouched Australian A-League rows changed too.
$ws.Range("G2").Value = 2.62
$ws.Range("L2").Value = 1.41
$ws.Range("N2").Value = 3.75
$ws.Range("W2").Value = 1.61

$ws.Range("L3").Value = 1.3
$ws.Range("R3").Value = 1.43

# 4) And a few odds changed on matches that were pushed down by the insert
#    (German Bundesliga -> now row 8, English National League -> now row 9,
#    Spanish La Liga -> now row 10).
$ws.Range("F8").Value = 1.53
$ws.Range("G8").Value = 1.55
$ws.Range("I8").Value = 7

$ws.Range("J9").Value = 3.6

$ws.Range("H10").Value = 4.2
$ws.Range("J10").Value = 3.5
$ws.Range("O10").Value = 1.47
$ws.Range("AI10").Value = 85
